$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple D/E-only updates (no swap rows)
$ws.Range("D2").Value = "29.439.94"

$ws.Range("D3").Value = "1.854.98"
$ws.Range("E3").Value = "  +1.14%  "

$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "244.82"
$ws.Range("E5").Value = "  -0.20%  "

$ws.Range("D6").Value = "0.6924"
$ws.Range("E6").Value = "  +0.18%  "

$ws.Range("D7").Value = "1.001"

$ws.Range("E8").Value = "  -0.32%  "

$ws.Range("E9").Value = "  +0.28%  "

$ws.Range("E10").Value = "  -0.07%  "

$ws.Range("D11").Value = "0.07765"
$ws.Range("E11").Value = "  -0.56%  "

# Row 12 and 13 swap (WrappedEther <-> Polkadot)
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "5.135"
$ws.Range("E12").Value = "  +1.11%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.847.73"
$ws.Range("E13").Value = "  +0.66%  "

$ws.Range("D14").Value = "0.6901"
$ws.Range("E14").Value = "  +1.40%  "

$ws.Range("D15").Value = "90.46"
$ws.Range("E15").Value = "  -0.19%  "

$ws.Range("D16").Value = "6.421"
$ws.Range("E16").Value = "  -0.31%  "

$ws.Range("D17").Value = "29.423.76"
$ws.Range("E17").Value = "  +1.69%  "

$ws.Range("D18").Value = "0.000008269"
$ws.Range("E18").Value = "  -1.01%  "

$ws.Range("D19").Value = "2.100.59"
$ws.Range("E19").Value = "  +0.87%  "

$ws.Range("D20").Value = "237.81"
$ws.Range("E20").Value = "  -2.27%  "

$ws.Range("E21").Value = "  +0.22%  "

$ws.Range("D22").Value = "1.001"

$ws.Range("D23").Value = "7.656"
$ws.Range("E23").Value = "  +2.40%  "

$ws.Range("E24").Value = "  +0.15%  "

$ws.Range("E25").Value = "  +1.40%  "

# Row 26 and 27 swap (Cosmos <-> Monero)
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "159.87"
$ws.Range("E26").Value = "  -1.39%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "8.889"
$ws.Range("E27").Value = "  +0.99%  "

$ws.Range("D28").Value = "18.25"
$ws.Range("E28").Value = "  +0.28%  "

$ws.Range("D29").Value = "1.535"
$ws.Range("E29").Value = "  -1.20%  "

$ws.Range("D30").Value = "4.248"
$ws.Range("E30").Value = "  +0.70%  "

$ws.Range("D31").Value = "4.152"
$ws.Range("E31").Value = "  -0.19%  "

$ws.Range("D32").Value = "1.193"
$ws.Range("E32").Value = "  +1.06%  "

$ws.Range("D33").Value = "0.05113"
$ws.Range("E33").Value = "  -0.32%  "

$ws.Range("D34").Value = "0.7683"
$ws.Range("E34").Value = "  +0.32%  "

$ws.Range("D35").Value = "1.887"
$ws.Range("E35").Value = "  +2.19%  "

$ws.Range("D36").Value = "1.149"
$ws.Range("E36").Value = "  +0.23%  "

$ws.Range("D37").Value = "2.686"
$ws.Range("E37").Value = "  +0.22%  "

$ws.Range("D38").Value = "1.330.09"
$ws.Range("E38").Value = "  +8.17%  "

$ws.Range("D39").Value = "0.01861"
$ws.Range("E39").Value = "  +0.98%  "

$ws.Range("D40").Value = "0.9751"
$ws.Range("E40").Value = "  +5.64%  "

$ws.Range("D41").Value = "2.712"
$ws.Range("E41").Value = "  +0.58%  "

$ws.Range("D42").Value = "106.07"
$ws.Range("E42").Value = "  -2.08%  "

$ws.Range("D43").Value = "5.840"
$ws.Range("E43").Value = "  -0.35%  "

$ws.Range("E44").Value = "  +0.18%  "

$ws.Range("E45").Value = "  +2.90%  "

$ws.Range("D46").Value = "9.765"
$ws.Range("E46").Value = "  +1.95%  "

$ws.Range("D47").Value = "2.002.88"
$ws.Range("E47").Value = "  +1.05%  "

$ws.Range("E48").Value = "  +0.96%  "

$ws.Range("D49").Value = "1.777"
$ws.Range("E49").Value = "  +1.81%  "

$ws.Range("D50").Value = "63.10"
$ws.Range("E50").Value = "  -1.27%  "

$ws.Range("D51").Value = "6.972"
$ws.Range("E51").Value = "  +0.56%  "
